$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data is inserted at the top of this sub-block (rows 211-212),
# pushing all the following weeks down by two rows (old row 321/322 lands on 323/324).
$ws.Rows("211:212").Insert()

# Row 211 - "Primera" quality, new week (fecha serial 44596)
$ws.Range("A211").Value = 8
$ws.Range("B211").Value = "Terminal La Palmera de La Serena"
$ws.Range("C211").Value = "Coquimbo"
$ws.Range("D211").Value = 44596
$ws.Range("E211").Value = 4
$ws.Range("F211").Value = 100112017
$ws.Range("G211").Value = "Apio"
$ws.Range("H211").Value = "Americana (o)"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 2200
$ws.Range("K211").Value = 8500
$ws.Range("L211").Value = 9000
$ws.Range("M211").Value = 8750
$ws.Range("N211").Value = "`$/docena de matas"
$ws.Range("O211").Value = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P211").Value = 1458
$ws.Range("Q211").Value = 6
$ws.Range("R211").Value = "Hortaliza"

# Row 212 - "Segunda" quality, same new week
$ws.Range("A212").Value = 8
$ws.Range("B212").Value = "Terminal La Palmera de La Serena"
$ws.Range("C212").Value = "Coquimbo"
$ws.Range("D212").Value = 44596
$ws.Range("E212").Value = 4
$ws.Range("F212").Value = 100112017
$ws.Range("G212").Value = "Apio"
$ws.Range("H212").Value = "Americana (o)"
$ws.Range("I212").Value = "Segunda"
$ws.Range("J212").Value = 1300
$ws.Range("K212").Value = 6500
$ws.Range("L212").Value = 7000
$ws.Range("M212").Value = 6750
$ws.Range("N212").Value = "`$/docena de matas"
$ws.Range("O212").Value = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P212").Value = 1125
$ws.Range("Q212").Value = 6
$ws.Range("R212").Value = "Hortaliza"
